# Updated Excel To-Do List with latest progress and status
#
# The workbook tracks a to-do list in "To Do List Check Box". Each task row
# (rows 5-29) has a Forms checkbox in column E that marks the task as done.
# Row 8 is the 4th task ("Check Box 5" / ctrlProp4.xml, linked to cell $E$8).
# Marking it complete flips E8 from FALSE to TRUE, which in turn feeds the
# existing worksheet formulas that total completed/weighted tasks:
#   H10 = COUNTIFS(...)                    -> 10  becomes 13
#   H12 = IFERROR(H10/H11,0)  (progress %) -> 10.9% becomes 14.1%

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tick the checkbox for task row 8 by setting its linked cell to TRUE.
$ws.Range("E8").Value = $True

# Keep the checkbox shape's own state in sync as well.
$checkBox = $ws.Shapes.Item("Check Box 5")
$checkBox.ControlFormat.Value = 1

# Make sure every dependent formula (counts, percentages) is up to date.
$excel.CalculateFullRebuild()

$wb.Save()
